# Results after investigating 528 ok pairs
#
# The "Summary" sheet's data table (A1:Q14) gets re-sorted (ascending) on
# columns D, E, G, H, J, K, M, N - mirroring a Data > Sort operation the
# author performed in Excel after adding rows 10-14 (the "sfnf" cross
# check pairs). The selection on that sheet moves to the newly produced
# P1:Q14 comments column, and the app window is resized taller.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Re-sort A1:Q14 using the same 8 ascending keys recorded in the sheet's
# sortState (D, E, G, H, J, K, M, N) - no header row in the sorted range.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("D1:D14")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("E1:E14")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("G1:G14")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("H1:H14")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("J1:J14")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("K1:K14")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("M1:M14")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("N1:N14")) | Out-Null
$ws.Sort.SetRange($ws.Range("A1:Q14"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Move the selection to the comment columns for the freshly sorted rows.
$ws.Activate()
$ws.Range("P1:Q14").Select()

# Resize the workbook's window (taller), matching the recorded
# bookViews/workbookView windowHeight change.
$win = $wb.Windows.Item(1)
$win.Height = 26980
$win.Width = 51200
